$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Profiling" row (row 3) entirely; rows below shift up.
$ws.Rows.Item(3).Delete()

# Update the Time column values to the new results.
$ws.Range("A2").Value = 3.54
$ws.Range("A3").Value = 11.52
$ws.Range("A4").Value = 129.3
$ws.Range("A5").Value = 45.36
